$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Coliflor at the
# "Macroferia Regional de Talca" market. It belongs chronologically
# right before the existing row 60 entry, so insert a fresh row there
# and push every subsequent record down by one (row 177 -> row 178).
$ws.Rows(60).Insert()

$row = 60
$ws.Cells.Item($row, 1).Value2  = 5
$ws.Cells.Item($row, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value2  = "Maule"
$ws.Cells.Item($row, 4).Value2  = 44533
$ws.Cells.Item($row, 5).Value2  = 7
$ws.Cells.Item($row, 6).Value2  = 100112008
$ws.Cells.Item($row, 7).Value2  = "Coliflor"
$ws.Cells.Item($row, 8).Value2  = "Sin especificar"
$ws.Cells.Item($row, 9).Value2  = "Primera"
$ws.Cells.Item($row, 10).Value2 = 3000
$ws.Cells.Item($row, 11).Value2 = 600
$ws.Cells.Item($row, 12).Value2 = 600
$ws.Cells.Item($row, 13).Value2 = 600
$ws.Cells.Item($row, 14).Value2 = "`$/unidad"
$ws.Cells.Item($row, 15).Value2 = "Región del Maule"
$ws.Cells.Item($row, 16).Value2 = 600
$ws.Cells.Item($row, 17).Value2 = 1
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"

# Match the date-formatted style used by every other row's "Fecha" cell.
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row + 1, 4).NumberFormat
